# EPBDS-9540 Support Java Name conversion on Json field name generating in SpreadsheetResults
# The step names used inside the _res_.$StepN["StepM"]:Integer markers (used to
# describe JSON field names for SpreadsheetResult fields) must be lower-cased
# (e.g. "Step1" -> "step1") to follow Java naming conventions for generated
# JSON field names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(37, 38, 56, 57)
$cols = @("C", "D", "E", "F", "G", "H")

foreach ($r in $rows) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $cell = $ws.Range($addr)
        $value = $cell.Value2
        if ($value -ne $null) {
            $newValue = $value -creplace '\["Step(\d+)"\]', '["step$1"]'
            if (-not $value.Equals($newValue)) {
                $cell.Value2 = $newValue
            }
        }
    }
}
